{"js": "const body = context.document.body;\nconst pairs = [\n  [\"Die Marketingkampagne wird die folgende Tagline verwenden, um Munsons Markenkern zu erfassen: \\\"Munsons: Pickles and Preserves with a Purpose\\\".\", \"Die Marketingkampagne wird die folgende Tagline verwenden, um Munsons Markenkern zu erfassen: \u201eMunson's: Pickles and Preserves with a Purpose\u201c.\"],\n  [\"Die Marketingkampagne wird den folgenden Slogan verwenden, um die Produktvorteile von Munson hervorzuheben: \\\"Munsons: Mehr als nur Pickles und Preserves\\\".\", \"Die Marketingkampagne wird den folgenden Slogan verwenden, um die Produktvorteile von Munson hervorzuheben: \u201eMunson's: More than Just Pickles and Preserves\u201c.\"],\n  [\"Die Marketingkampagne wird das folgende Motto verwenden, um Munsons Kundenvertretung zu inspirieren: \\\"Munsons: Teilen der Liebe von Pickles und Preserves\\\".\", \"Die Marketingkampagne wird das folgende Motto verwenden, um die Customer Advocacy von Munson's zu unterst\u00fctzen: \u201eMunson's: Share the Love of Pickles and Preserves\u201c.\"],\n  [\"Die Marketingkampagne wird den folgenden Ausdruck verwenden, um Munsons Produktabonnement zu testen und zu kaufen: \\\"Munsons: Find Them, Try Them, Love Them\\\".\", \"Die Marketingkampagne wird die folgende Phrase verwenden, um das Produktabonnement von Munson's zu testen und zu kaufen: \u201eMunson's: Find Them, Try Them, Love Them\u201c.\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    @('Die Marketingkampagne wird die folgende Tagline verwenden, um Munsons Markenkern zu erfassen: \"Munsons: Pickles and Preserves with a Purpose\".', 'Die Marketingkampagne wird die folgende Tagline verwenden, um Munsons Markenkern zu erfassen: \u201eMunson''s: Pickles and Preserves with a Purpose\u201c.'),\n    @('Die Marketingkampagne wird den folgenden Slogan verwenden, um die Produktvorteile von Munson hervorzuheben: \"Munsons: Mehr als nur Pickles und Preserves\".', 'Die Marketingkampagne wird den folgenden Slogan verwenden, um die Produktvorteile von Munson hervorzuheben: \u201eMunson''s: More than Just Pickles and Preserves\u201c.'),\n    @('Die Marketingkampagne wird das folgende Motto verwenden, um Munsons Kundenvertretung zu inspirieren: \"Munsons: Teilen der Liebe von Pickles und Preserves\".', 'Die Marketingkampagne wird das folgende Motto verwenden, um die Customer Advocacy von Munson''s zu unterst\u00fctzen: \u201eMunson''s: Share the Love of Pickles and Preserves\u201c.'),\n    @('Die Marketingkampagne wird den folgenden Ausdruck verwenden, um Munsons Produktabonnement zu testen und zu kaufen: \"Munsons: Find Them, Try Them, Love Them\".', 'Die Marketingkampagne wird die folgende Phrase verwenden, um das Produktabonnement von Munson''s zu testen und zu kaufen: \u201eMunson''s: Find Them, Try Them, Love Them\u201c.'),\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $found = $find.Execute()\n    if ($found) {\n        $range.Text = $newText\n    } else {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n"}
